$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Il1a -> Il1r2, Target cluster ECs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02286966666666667
$ws.Range("H2").Value = 0.068609
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.005328
$ws.Range("N2").Value = 0.015984
$ws.Range("O2").Value = 0.001689940172269439
$ws.Range("P2").Value = 0.001689940172269439
$ws.Range("Q2").Value = 0.000121849584
$ws.Range("R2").Value = 0.001096646256
$ws.Range("S2").Value = 0.001689940172269439
$ws.Range("T2").Value = 0.001689940172269439

# Row 3 updates (Il1a -> Il1r2, Target cluster FAPs)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02286966666666667
$ws.Range("H3").Value = 0.068609
$ws.Range("O3").Value = 0.9983100598277306
$ws.Range("P3").Value = 0.9983100598277306
$ws.Range("Q3").Value = 0.07198104849455556
$ws.Range("R3").Value = 0.6478294364510001
$ws.Range("S3").Value = 0.9983100598277306
$ws.Range("T3").Value = 0.9983100598277306
